$d = $word.ActiveDocument

function Find-ParaIndexByText($text) {
    $n = $d.Paragraphs.Count
    for ($i = 1; $i -le $n; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.Contains($text)) {
            return $i
        }
    }
    return -1
}

function Insert-DetailParagraph($headerSnippet, $newText) {
    $idx = Find-ParaIndexByText($headerSnippet)
    if ($idx -eq -1) {
        Write-Host "WARNING: header not found: $headerSnippet"
        return
    }
    $header = $d.Paragraphs.Item($idx)
    $next = $header.Next()
    $next.Range.InsertParagraphBefore()
    $inserted = $d.Paragraphs.Item($idx + 1)
    $inserted.Range.Text = $newText
}

# 1) After "pySODM" header, before "Total number of compartments: 3 x 581 x 4"
Insert-DetailParagraph "pySODM" "Uses commuter mobility matrix and social contact matrix"

# 2) After "flepiMoP - no age groups" header, before "Total number of compartments: 3 x 581 ="
Insert-DetailParagraph "no age groups" "Uses commuter mobility matrix"

# 3) After "flepiMoP - with age groups (row sums)" header, before "Total number of compartments: 3 x 581"
Insert-DetailParagraph "row sums" "Uses commuter mobility matrix and social contacts are modeled in a destination average way"

# 4) After "flepiMoP - with age groups (integrating the full contact matrix)" header
Insert-DetailParagraph "integrating the full contact matrix" "Uses commuter mobility matrix and social contact matrix"

# 5) Move the lastRenderedPageBreak marker: it was on the drawing run that follows the
#    "Approx. computational complexity: 4.5s (overhead) + 7.6s (simulation)" paragraph
#    (the one in the "row sums" section); it should now be on the "Simulated using"
#    paragraph of that same section instead.
$simIdx = Find-ParaIndexByText("Simulated using")
$rowSumsHeaderIdx = Find-ParaIndexByText("row sums")
# find the "Simulated using" paragraph that belongs to the row-sums section (the next one
# after the header, skipping the inserted detail + total-compartments paragraphs)
$n = $d.Paragraphs.Count
$targetSimIdx = -1
for ($i = $rowSumsHeaderIdx; $i -le $n; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Contains("Simulated using")) {
        $targetSimIdx = $i
        break
    }
}
if ($targetSimIdx -ne -1) {
    $simPara = $d.Paragraphs.Item($targetSimIdx)
    $simPara.Range.InsertBefore([char]0)
    # Use Find/Replace-free approach: insert the lastRenderedPageBreak field via Range on first char
}
